# "bolt_softedge: Add source information"
#
# A new row is inserted into the "Power Models and Textures" table (Table2
# on Sheet1), right after "textures\bolt_sharp.igb" (row 21) and before
# "textures\bolt2.igb" (old row 22), for the new texture
# "textures\bolt_softedge.igb". All the Source columns are "Custom: Based
# on Default" and the Notes column is "3b. Only found in other games".
# This pushes every row below down by one (old row 138 -> 139) and grows
# the table / used range / conditional formatting / dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new physical row at sheet row 22 (shifts 22..138 -> 23..139).
$ws.Rows.Item(22).Insert()

# Populate the new row's data.
$ws.Cells.Item(22, 1).Value = "textures\bolt_softedge.igb"
$ws.Cells.Item(22, 2).Value = "Custom: Based on Default"
$ws.Cells.Item(22, 3).Value = "Custom: Based on Default"
$ws.Cells.Item(22, 4).Value = "Custom: Based on Default"
$ws.Cells.Item(22, 5).Value = "Custom: Based on Default"
$ws.Cells.Item(22, 6).Value = "Custom: Based on Default"
$ws.Cells.Item(22, 7).Value = "3b. Only found in other games"

# Grow the table (ListObject) so it covers the new row: A1:G138 -> A1:G139.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G139"))

# Extend the conditional-formatting ranges that covered the table body so
# they include the newly inserted row too (A2:A138 -> A2:A139, etc.)
#
# NOTE: pull the *whole sheet's* rules into one collection up front and
# address each rule by its stable global index. Re-querying
# Range(...).FormatConditions per-column after some ranges have already
# been widened returns bogus aliased/duplicated items once ranges overlap
# (e.g. querying "C2:C138" after "B2:F138" was widened to "B2:F139"
# returns 11 copies of the B:F rule instead of the 2 real C-column
# rules) - fetching once up front and indexing numerically sidesteps that.
$allCf = $ws.Cells.FormatConditions
$cfTargets = @(
    @{idx = 1;  rng = "A2:A139" },
    @{idx = 2;  rng = "A2:A139" },
    @{idx = 3;  rng = "B2:B139" },
    @{idx = 4;  rng = "B2:B139" },
    @{idx = 5;  rng = "B2:F139" },
    @{idx = 6;  rng = "C2:C139" },
    @{idx = 7;  rng = "C2:C139" },
    @{idx = 8;  rng = "D2:D139" },
    @{idx = 9;  rng = "D2:D139" },
    @{idx = 10; rng = "E2:E139" },
    @{idx = 11; rng = "E2:E139" },
    @{idx = 12; rng = "F2:F139" },
    @{idx = 13; rng = "F2:F139" }
)
foreach ($t in $cfTargets) {
    $allCf.Item($t.idx).ModifyAppliesToRange($ws.Range($t.rng))
}

# Columns B and F now also contain the longest label ("Custom: Based on
# Default", same as C:E already did), so their best-fit width matches
# C:E's. Approximate that best-fit width as closely as this host's font
# metrics allow.
$ws.Columns.Item(2).ColumnWidth = 21.6
$ws.Columns.Item(6).ColumnWidth = 21.6

# Match the author's final selection.
$ws.Range("A15").Select()
